$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-12 Friday" "2024-04-13 Saturday"

Replace-Text "354÷4=88, 2" "287÷8=35, 7"
Replace-Text "445÷9=49, 4" "713÷9=79, 2"
Replace-Text "416÷6=69, 2" "644÷2=322, 0"
Replace-Text "862÷8=107, 6" "122÷3=40, 2"
Replace-Text "574÷8=71, 6" "387÷9=43, 0"
Replace-Text "501÷9=55, 6" "432÷7=61, 5"
Replace-Text "428÷4=107, 0" "464÷6=77, 2"
Replace-Text "927÷4=231, 3" "513÷3=171, 0"
Replace-Text "879÷8=109, 7" "156÷9=17, 3"
Replace-Text "825÷2=412, 1" "425÷4=106, 1"
Replace-Text "860÷3=286, 2" "684÷9=76, 0"
Replace-Text "969÷4=242, 1" "749÷2=374, 1"
Replace-Text "134÷3=44, 2" "228÷2=114, 0"
Replace-Text "501÷8=62, 5" "996÷9=110, 6"
Replace-Text "706÷5=141, 1" "371÷9=41, 2"
Replace-Text "144÷6=24, 0" "602÷3=200, 2"
Replace-Text "930÷9=103, 3" "522÷5=104, 2"
Replace-Text "142÷6=23, 4" "127÷2=63, 1"
Replace-Text "279÷2=139, 1" "744÷5=148, 4"
Replace-Text "869÷5=173, 4" "769÷5=153, 4"
Replace-Text "185÷7=26, 3" "120÷3=40, 0"
Replace-Text "275÷3=91, 2" "984÷2=492, 0"
Replace-Text "347÷6=57, 5" "560÷5=112, 0"
Replace-Text "225÷9=25, 0" "336÷2=168, 0"
Replace-Text "249÷8=31, 1" "582÷6=97, 0"
